$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 266; this shifts the existing rows 266-283 down to 267-284
$ws.Rows.Item(266).Insert()

# Populate the newly inserted row 266 with the new weekly data point
$ws.Cells.Item(266, 1).Value2 = 4
$ws.Cells.Item(266, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(266, 3).Value = "Los Lagos"
$ws.Cells.Item(266, 4).Value2 = 44783
$ws.Cells.Item(266, 5).Value2 = 10
$ws.Cells.Item(266, 6).Value = "Fruta"
$ws.Cells.Item(266, 7).Value2 = 100108
$ws.Cells.Item(266, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(266, 9).Value2 = 100108005
$ws.Cells.Item(266, 10).Value = "Piña"
$ws.Cells.Item(266, 11).Value = "Caramelo"
$ws.Cells.Item(266, 12).Value = "Primera"
$ws.Cells.Item(266, 13).Value2 = 60
$ws.Cells.Item(266, 14).Value2 = 23000
$ws.Cells.Item(266, 15).Value2 = 23000
$ws.Cells.Item(266, 16).Value2 = 23000
$ws.Cells.Item(266, 17).Value = "$/caja 12 unidades"
$ws.Cells.Item(266, 18).Value = "Ecuador"
$ws.Cells.Item(266, 19).Value2 = 1917
$ws.Cells.Item(266, 20).Value2 = 12
